# The presentation's design was changed away from the green "Integral"
# theme back to the default "Office Theme" colour palette.
#
# ppt/theme/theme2.xml is the theme part that is actually wired up to the
# slide master (and to the presentation's own theme relationship), so it is
# the file that must end up holding the "Office Theme" colours. This COM
# shim does not give us file-system access to load a .thmx, nor does it
# persist ThemeColorScheme.Name / Theme.Name edits, so the supported way to
# rewrite a theme in this environment is exactly what's documented on the
# Theme object: ThemeColorScheme.Colors(i).RGB (and ThemeFontScheme
# MajorFont/MinorFont for fonts, unneeded here since both themes already
# share the same font scheme).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$scheme = $theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB() long values (0x00BBGGRR).
$scheme.Colors(1).RGB  = 0        # dk1       000000
$scheme.Colors(2).RGB  = 16777215 # lt1       FFFFFF
$scheme.Colors(3).RGB  = 6968388  # dk2       44546A
$scheme.Colors(4).RGB  = 15132391 # lt2       E7E6E6
$scheme.Colors(5).RGB  = 13998939 # accent1   5B9BD5
$scheme.Colors(6).RGB  = 3243501  # accent2   ED7D31
$scheme.Colors(7).RGB  = 10855845 # accent3   A5A5A5
$scheme.Colors(8).RGB  = 49407    # accent4   FFC000
$scheme.Colors(9).RGB  = 12874308 # accent5   4472C4
$scheme.Colors(10).RGB = 4697456  # accent6   70AD47
$scheme.Colors(11).RGB = 12673797 # hlink     0563C1
$scheme.Colors(12).RGB = 7491477  # folHlink  954F72
